$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-Text "2025-02-25 Tuesday" "2025-02-26 Wednesday"

Replace-Text "837×8=6696" "149×8=1192"
Replace-Text "238×7=1666" "358×3=1074"
Replace-Text "343×5=1715" "653×2=1306"
Replace-Text "459×5=2295" "411×2=822"
Replace-Text "516×9=4644" "344×7=2408"
Replace-Text "776×8=6208" "622×9=5598"
Replace-Text "398×9=3582" "115×3=345"
Replace-Text "468×2=936" "289×6=1734"
Replace-Text "538×5=2690" "162×4=648"
Replace-Text "192×8=1536" "931×3=2793"
Replace-Text "707×7=4949" "297×6=1782"
Replace-Text "269×9=2421" "331×7=2317"
Replace-Text "472×3=1416" "674×4=2696"
Replace-Text "881×9=7929" "854×6=5124"
Replace-Text "151×5=755" "542×3=1626"
Replace-Text "569×4=2276" "174×8=1392"
Replace-Text "174×2=348" "154×4=616"
Replace-Text "195×5=975" "655×3=1965"
Replace-Text "519×3=1557" "911×4=3644"
Replace-Text "872×4=3488" "831×9=7479"
Replace-Text "193×4=772" "788×9=7092"
Replace-Text "965×3=2895" "314×3=942"
Replace-Text "358×4=1432" "794×3=2382"
Replace-Text "966×5=4830" "792×4=3168"
Replace-Text "889×3=2667" "629×3=1887"
